$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header columns (P, Q, R) to support importing values for additional labels
$ws.Range("P1").Value = "Label Rouge"
$ws.Range("Q1").Value = "AOC / AOP / IGP"
$ws.Range("R1").Value = "HVE"

# Fill in example values for the existing example rows
$ws.Range("P2").Value = 20
$ws.Range("P3").Value = 30

# Match the column width that Excel auto-calculated for the new "AOC / AOP / IGP" column
$ws.Columns.Item(17).ColumnWidth = 13.8

# Leave the selection on the last edited cell, as in the source edit
$ws.Range("P3").Select() | Out-Null
